# "added: ppd lab 5"
# Update the three service/location values on the "Input" sheet and move the
# active selection, matching the author's edits to InputFile.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Row 3 - service type changed from "Tuns" (haircut) to "Barba" (beard)
$ws.Range("C3").Value = "Barba"

# Row 7 - service type changed from "Masaj" to "Masaj barbati" (men's massage)
$ws.Range("C7").Value = "Masaj barbati"

# Row 9 - location renamed from "Targu Mures" to "Targu-Mures"
$ws.Range("B9").Value = "Targu-Mures"

# Move the active selection on the Input sheet to B23, matching the saved
# cursor position in the workbook.
$ws.Activate() | Out-Null
$ws.Range("B23").Select() | Out-Null

Write-Output "done"
